# Fix the "transformer" testcase sheet: it was missing the "b" (susceptance)
# column that every other sheet already accounts for. Insert a new column H
# with header "b", and restore the originally-intended active sheet/selection
# state (transformer tab active with H1 selected; generator tab left on C54,
# no longer the active tab).

$wb = $excel.ActiveWorkbook

# First move the selection on the "generator" sheet to C54 (its new resting
# selection) without leaving it the active tab.
$wsGenerator = $wb.Worksheets.Item("generator")
$wsGenerator.Range("C54").Select()

# Now activate "transformer", insert the missing "b" column after the
# existing "x" column (before the old H1/"ShortTermRating" column), and
# select the new header cell.
$wsTransformer = $wb.Worksheets.Item("transformer")
$wsTransformer.Activate()
$wsTransformer.Range("H1").EntireColumn.Insert()
$wsTransformer.Range("H1").Value = "b"
$wsTransformer.Range("H1").Select()
